$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.059568881988525
$ws.Range("B1").Value = 2.846658945083618
$ws.Range("C1").Value = 5.189640522003174
$ws.Range("D1").Value = 3.583705902099609
$ws.Range("E1").Value = 1.314777612686157
